# Sincronizando o repositório local baixado do Luciano com o repositório que eu criei
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=22223; B="Raquel Aragão";        C="Financeiro";             D="Doença";             E=3; F=45082; G=8831.67}
    @{Row=3;  A=69056; B="Maria Vitória Farias";  C="TI";                     D="Doença";             E=1; F=45105; G=9152.48}
    @{Row=4;  A=66129; B="Benjamin Carvalho";     C="Engenharia";             D="Problemas pessoais"; E=8; F=45093; G=10814.64}
    @{Row=5;  A=58602; B="Erick Moraes";          C="Jurídico";               D="Outros";             E=3; F=45100; G=4267.34}
    @{Row=6;  A=12254; B="Gabrielly Vieira";      C="Recursos Humanos";       D="Viagem de negócios"; E=6; F=45105; G=11732.68}
    @{Row=7;  A=72722; B="Pietra da Rosa";        C="Engenharia";             D="Doença";             E=8; F=45078; G=9097.76}
    @{Row=8;  A=75883; B="Maria Sophia Pinto";    C="Jurídico";               D="Consulta médica";    E=1; F=45090; G=8745.42}
    @{Row=9;  A=66421; B="Mariana Mendes";        C="Financeiro";             D="Outros";             E=1; F=45102; G=6080.97}
    @{Row=10; A=43747; B="Anthony da Cunha";      C="P&D";                    D="Problemas pessoais"; E=1; F=45088; G=11400.88}
    @{Row=11; A=87757; B="Luiza Porto";           C="Atendimento ao Cliente"; D="Doença";             E=1; F=45091; G=8350.280000000001}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
